# Update the "pubmed" reference column (K) on the VOX sheet: the publication
# labels are reformatted from space-separated to underscore-separated text.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VOX")

$ws.Range("K2:K3").Value = "AASLD_2016_Abs_846"
$ws.Range("K4").Value = "Pooled_phase_III_EASL_2017_Abs_248"
$ws.Range("K5:K102").Value = "AASLD_2017_Abs_1176"
$ws.Range("K103:K107").Value = "EASL_2017_Abs_THU-257"
$ws.Range("K108:K117").Value = "AASLD_2015_Abs_718"

# Scroll the view so that G92 is the top-left visible cell, with G92 active
# and the whole sheet selected (mirrors the saved view state in the workbook).
$ws.Activate()
$ws.Range("A1:XFD1048576").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 92
$excel.ActiveWindow.ScrollColumn = 7
$ws.Range("G92").Activate() | Out-Null
